$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 (pushes old rows 7..34 down to 8..35)
$ws.Rows(7).Insert()

# Populate the newly inserted row 7. It is the same record as the
# (now shifted) row 8 below it, except for a new sampling date (D) and
# a different Volumen (M).
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Femacal de La Calera"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44453
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108004
$ws.Range("J7").Value = "Papaya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia del Elquí"
$ws.Range("S7").Value = 1500
$ws.Range("T7").Value = 10
